$d = $word.ActiveDocument

# 1. "The brute force was done allowing" -> "The brute force was done by allowing"
$r = $d.Content
$r.Find.Execute("The brute force was done allowing", $true, $false, $false, $false, $false, $true, 1, $false, "The brute force was done by allowing", 2)

# 2. "wrong password, then we estimated" -> "wrong password, and then we estimated"
$r = $d.Content
$r.Find.Execute("wrong password, then we estimated", $true, $false, $false, $false, $false, $true, 1, $false, "wrong password, and then we estimated", 2)

# 3. "calculating how long it would take to try" -> "calculating the time it would take to try"
$r = $d.Content
$r.Find.Execute("calculating how long it would take to try", $true, $false, $false, $false, $false, $true, 1, $false, "calculating the time it would take to try", 2)

# 4. Move the "_GoBack" bookmark to reflect the new last-edit location: right after
#    "Based on average time and stan" in the Conclusion paragraph (mid "standard").
$r = $d.Content
$r.Find.Execute("Based on average time and stan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
